# Add a new "Data_Source" column (M) to the master summary sheet,
# marking every existing data row as sourced from the "Network".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: copy the formatting used by the other header cells (bold,
# centered, bordered) from A1 onto the new M1 header, then set its text.
$ws.Range("A1").Copy()
$ws.Range("M1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("M1").Value = "Data_Source"

# Data rows: tag each existing row as coming from the "Network" source.
$ws.Range("M2").Value = "Network"
$ws.Range("M3").Value = "Network"
$ws.Range("M4").Value = "Network"
$ws.Range("M5").Value = "Network"

$excel.CutCopyMode = 0
